# ASW1 Addressing opdateret for LAB 4.3.8
#
# Fills in the ASW2 port/VLAN table and nudges the DSW1 sheet's saved
# selection, matching the upstream commit.

$wb = $excel.ActiveWorkbook

# --- ASW2: fill in the switch port / VLAN addressing table -----------------
$ws = $wb.Worksheets.Item("ASW2")
$ws.Activate()

$ws.Range("A2").Value = "ASW2"
$ws.Range("E2").Value = "192.168.1.130"

$ws.Range("B3").Value = "Gi1/0/1-2(23)"
$ws.Range("C3").Value = "VLAN 20 Access Mode"

$ws.Range("B4").Value = "Gi1/0/19"
$ws.Range("C4").Value = "TRUNK Native 99"

$ws.Range("B5").Value = "Gi1/0/21"
$ws.Range("C5").Value = "TRUNK Native 99"

$ws.Range("B6").Value = "Gi1/0/23"
$ws.Range("C6").Value = "Access20 Voice150"

$ws.Range("B7").Value = "Gi1/0/24"
$ws.Range("C7").Value = "VLAN 250 Access Mode"

$ws.Range("B8").Value = "VLAN 10"
$ws.Range("C8").Value = "N/A"

$ws.Range("B9").Value = "VLAN 20"
$ws.Range("C9").Value = "N/A"

$ws.Range("B10").Value = "VLAN 150"
$ws.Range("C10").Value = "N/A"

$ws.Range("B11").Value = "VLAN 200"
$ws.Range("C11").Value = "192.168.1.132"
$ws.Range("D11").Value = 255255255224
$ws.Range("D11").NumberFormat = "#,##0"

$ws.Range("B12").Value = "VLAN 250"
$ws.Range("C12").Value = "N/A"

# Widen column C to fit the new text (Excel's "best fit" double-click-the-
# border autosize), and leave the saved selection where the editing
# session ended up.
$ws.Columns("C").ColumnWidth = 20.6
$ws.Range("C8").Select()

# --- DSW1: selection moved while reviewing ---------------------------------
$ws2 = $wb.Worksheets.Item("DSW1")
$ws2.Activate()
$ws2.Range("D14").Select()

# Leave ASW2 as the active/visible sheet, matching tabSelected="1" there.
$ws.Activate()
